$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Insert a new row at position 11 ("Jurisdiction" / empty), pushing the
# existing rows (Description, Purpose, Copyright, ...) down by one.
$ws.Rows.Item(11).Insert()

# The freshly inserted row picks up a generic style; copy the same
# formatting used by the surrounding metadata rows (row 12, "Description").
$ws.Range("A12:B12").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)

# Populate the new row.
$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = ""

# Refresh the Date metadata value.
$ws.Range("B8").Value = "2024-07-01T07:50:29+00:00"
